# DB Change for Linux
$wb = $excel.ActiveWorkbook

# --- addListItem sheet: update source value LinuxAA -> LinuxAB ---
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "LinuxAB"

# --- createUser sheet: bump the numeric id used to build the test user 2710 -> 2711 ---
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 2711

# Recalculate so cached formula results (CONCAT/A2 references) are refreshed
$excel.Calculate()

# --- Selections / active sheet ---
# loginDetails: selection moves to B2
$wsLogin = $wb.Worksheets.Item("loginDetails")
$wsLogin.Activate()
$wsLogin.Range("B2").Select()

# addListItem: selection stays at A2, but it is no longer the tab-selected sheet
$wsAddListItem.Activate()
$wsAddListItem.Range("A2").Select()

# createUser: becomes the tab-selected / active sheet, selection at A2
$wsCreateUser.Activate()
$wsCreateUser.Range("A2").Select()
